$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize Spanish connector words (de, del, la, el, los, las, y) in municipality/state names
$renames = @(
    @('B8', 'Pabellón De Arteaga'),
    @('B9', 'Rincón De Romos'),
    @('B10', 'San Francisco De Los Romo'),
    @('B11', 'San José De Gracia'),
    @('B16', 'Playas De Rosarito'),
    @('B42', 'Amatenango De La Frontera'),
    @('B43', 'Amatenango Del Valle'),
    @('B46', 'Bejucal De Ocampo'),
    @('B48', 'Benemérito De Las Américas'),
    @('B58', 'Chiapa De Corzo'),
    @('B65', 'Comitán De Domínguez'),
    @('B94', 'Marqués De Comillas'),
    @('B95', 'Mazapa De Madero'),
    @('B98', 'Montecristo De Guerrero'),
    @('B103', 'Ocozocoautla De Espinosa'),
    @('B115', 'Salto De Agua'),
    @('B116', 'San Cristóbal De Las Casas'),
    @('B163', 'Coyame Del Sotol'),
    @('B174', 'Guadalupe Y Calvo'),
    @('B177', 'Hidalgo Del Parral'),
    @('B201', 'San Francisco De Borja'),
    @('B202', 'San Francisco De Conchos'),
    @('B203', 'San Francisco Del Oro'),
    @('B211', 'Valle De Zaragoza'),
    @('B236', 'San Juan De Sabinas'),
    @('B252', 'Villa De Álvarez'),
    @('A254', 'Ciudad De México'),
    @('B258', 'Cuajimalpa De Morelos'),
    @('B273', 'Coneto De Comonfort'),
    @('B287', 'Nombre De Dios'),
    @('B291', 'Pánuco De Coronado'),
    @('B298', 'San Juan De Guadalupe'),
    @('B299', 'San Juan Del Río'),
    @('B300', 'San Luis Del Cordero'),
    @('B301', 'San Pedro Del Gallo'),
    @('A311', 'Estado De México'),
    @('B311', 'Acambay De Ruíz Castañeda'),
    @('B314', 'Almoloya De Alquisiras'),
    @('B315', 'Almoloya De Juárez'),
    @('B316', 'Almoloya Del Río'),
    @('B323', 'Atizapán De Zaragoza'),
    @('B331', 'Chapa De Mota'),
    @('B337', 'Coacalco De Berriozábal'),
    @('B344', 'Ecatepec De Morelos'),
    @('B352', 'Ixtapan De La Sal'),
    @('B353', 'Ixtapan Del Oro'),
    @('B370', 'Naucalpan De Juárez'),
    @('B384', 'San Antonio La Isla'),
    @('B385', 'San Felipe Del Progreso'),
    @('B386', 'San Martín De Las Pirámides'),
    @('B388', 'San Simón De Guerrero'),
    @('B390', 'Soyaniquilpan De Juárez'),
    @('B400', 'Tenango Del Aire'),
    @('B401', 'Tenango Del Valle'),
    @('B414', 'Tlalnepantla De Baz'),
    @('B420', 'Valle De Bravo'),
    @('B421', 'Valle De Chalco Solidaridad'),
    @('B422', 'Villa De Allende'),
    @('B423', 'Villa Del Carbón'),
    @('B436', 'San Miguel De Allende'),
    @('B437', 'Apaseo El Alto'),
    @('B438', 'Apaseo El Grande'),
    @('B446', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B450', 'Jaral Del Progreso'),
    @('B458', 'Purísima Del Rincón'),
    @('B462', 'San Diego De La Unión'),
    @('B464', 'San Francisco Del Rincón'),
    @('B466', 'San Luis De La Paz'),
    @('B468', 'Santa Cruz De Juventino Rosas'),
    @('B470', 'Silao De La Victoria'),
    @('B475', 'Valle De Santiago'),
    @('B481', 'Acapulco De Juárez'),
    @('B484', 'Ajuchitlán Del Progreso'),
    @('B485', 'Alcozauca De Guerrero'),
    @('B489', 'Atenango Del Río'),
    @('B490', 'Atlamajalcingo Del Monte'),
    @('B492', 'Atoyac De Álvarez'),
    @('B493', 'Ayutla De Los Libres'),
    @('B496', 'Buenavista De Cuéllar'),
    @('B497', 'Chilapa De Álvarez'),
    @('B498', 'Chilpancingo De Los Bravo'),
    @('B499', 'Coahuayutla De José María Izazaga'),
    @('B504', 'Coyuca De Benítez'),
    @('B505', 'Coyuca De Catalán'),
    @('B509', 'Cuetzala Del Progreso'),
    @('B510', 'Cutzamala De Pinzón'),
    @('B516', 'Huitzuco De Los Figueroa'),
    @('B517', 'Iguala De La Independencia'),
    @('B519', 'Ixcateopan De Cuauhtémoc'),
    @('B520', 'Zihuatanejo De Azueta'),
    @('B522', 'La Unión De Isidoro Montes De Oca'),
    @('B525', 'Mártir De Cuilapan'),
    @('B538', 'Taxco De Alarcón'),
    @('B540', 'Técpan De Galeana'),
    @('B542', 'Tepecoacuilco De Trujano'),
    @('B544', 'Tixtla De Guerrero'),
    @('B548', 'Tlalixtaquilla De Maldonado'),
    @('B549', 'Tlapa De Comonfort'),
    @('B561', 'Agua Blanca De Iturbide'),
    @('B568', 'Atotonilco De Tula'),
    @('B569', 'Atotonilco El Grande'),
    @('B575', 'Cuautepec De Hinojosa'),
    @('B581', 'Huasca De Ocampo'),
    @('B585', 'Huejutla De Reyes'),
    @('B588', 'Jacala De Ledezma'),
    @('B595', 'Mineral De La Reforma'),
    @('B596', 'Mineral Del Chico'),
    @('B597', 'Mineral Del Monte'),
    @('B598', 'Mixquiahuala De Juárez'),
    @('B599', 'Molango De Escamilla'),
    @('B601', 'Nopala De Villagrán'),
    @('B602', 'Omitlán De Juárez'),
    @('B603', 'Pachuca De Soto'),
    @('B606', 'Progreso De Obregón'),
    @('B612', 'Santiago De Anaya'),
    @('B613', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B617', 'Tenango De Doria'),
    @('B619', 'Tepehuacán De Guerrero'),
    @('B620', 'Tepeji Del Río De Ocampo'),
    @('B623', 'Tezontepec De Aldama'),
    @('B632', 'Tula De Allende'),
    @('B633', 'Tulancingo De Bravo'),
    @('B634', 'Villa De Tezontepec'),
    @('B638', 'Zacualtipán De Ángeles'),
    @('B639', 'Zapotlán De Juárez'),
    @('B644', 'Acatlán De Juárez'),
    @('B645', 'Ahualulco De Mercado'),
    @('B650', 'Atemajac De Brizuela'),
    @('B653', 'Atotonilco El Alto'),
    @('B655', 'Autlán De Navarro'),
    @('B661', 'Cañadas De Obregón'),
    @('B668', 'Concepción De Buenos Aires'),
    @('B669', 'Cuautitlán De García Barragán'),
    @('B678', 'Encarnación De Díaz'),
    @('B685', 'Huejuquilla El Alto'),
    @('B686', 'Ixtlahuacán De Los Membrillos'),
    @('B687', 'Ixtlahuacán Del Río'),
    @('B691', 'Jilotlán De Los Dolores'),
    @('B697', 'La Manzanilla De La Paz'),
    @('B698', 'Lagos De Moreno'),
    @('B706', 'Ojuelos De Jalisco'),
    @('B711', 'San Cristóbal De La Barranca'),
    @('B712', 'San Diego De Alejandría'),
    @('B714', 'San Juan De Los Lagos'),
    @('B715', 'San Juanito De Escobedo'),
    @('B718', 'San Martín De Bolaños'),
    @('B720', 'San Miguel El Alto'),
    @('B721', 'San Sebastián Del Oeste'),
    @('B722', 'Santa María De Los Ángeles'),
    @('B723', 'Santa María Del Oro'),
    @('B726', 'Talpa De Allende'),
    @('B727', 'Tamazula De Gordiano'),
    @('B730', 'Techaluta De Montenegro'),
    @('B734', 'Teocuitatlán De Corona'),
    @('B735', 'Tepatitlán De Morelos'),
    @('B738', 'Tizapán El Alto'),
    @('B739', 'Tlajomulco De Zúñiga'),
    @('B751', 'Unión De San Antonio'),
    @('B752', 'Unión De Tula'),
    @('B753', 'Valle De Guadalupe'),
    @('B754', 'Valle De Juárez'),
    @('B759', 'Yahualica De González Gallo'),
    @('B760', 'Zacoalco De Torres'),
    @('B763', 'Zapotitlán De Vadillo'),
    @('B764', 'Zapotlán Del Rey'),
    @('B765', 'Zapotlán El Grande'),
    @('B791', 'Coalcomán De Vázquez Pallares'),
    @('B793', 'Cojumatlán De Régules'),
    @('B860', 'Tiquicheo De Nicolás Romero'),
    @('B886', 'Coatlán Del Río'),
    @('B894', 'Jonacatepec De Leandro Valle'),
    @('B898', 'Puente De Ixtla'),
    @('B904', 'Tetela Del Volcán'),
    @('B906', 'Tlaltizapán De Zapata'),
    @('B914', 'Zacualpan De Amilpas'),
    @('B918', 'Amatlán De Cañas'),
    @('B919', 'Bahía De Banderas'),
    @('B923', 'Ixtlán Del Río'),
    @('B930', 'Santa María Del Oro'),
    @('B962', 'San Nicolás De Los Garza'),
    @('B967', 'Acatlán De Pérez Figueroa'),
    @('B975', 'Ayoquezco De Aldama'),
    @('B979', 'Capulálpam De Méndez'),
    @('B981', 'Chalcatongo De Hidalgo'),
    @('B982', 'Chiquihuitlán De Benito Juárez'),
    @('B983', 'Ciénega De Zimatlán'),
    @('B986', 'Coicoyán De Las Flores'),
    @('B989', 'Constancia Del Rosario'),
    @('B992', 'Cuilápam De Guerrero'),
    @('B993', 'Cuyamecalco Villa De Zaragoza'),
    @('B994', 'El Barrio De La Soledad'),
    @('B996', 'Eloxochitlán De Flores Magón'),
    @('B997', 'Fresnillo De Trujano'),
    @('B998', 'Guadalupe De Ramírez'),
    @('B1000', 'Guelatao De Juárez'),
    @('B1001', 'Guevea De Humboldt'),
    @('B1002', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B1003', 'Heroica Ciudad De Huajuapan De León'),
    @('B1004', 'Heroica Ciudad De Tlaxiaco'),
    @('B1005', 'Huautla De Jiménez'),
    @('B1007', 'Ixtlán De Juárez'),
    @('B1008', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B1022', 'Magdalena Yodocono De Porfirio Díaz'),
    @('B1024', 'Mariscala De Juárez'),
    @('B1025', 'Mártires De Tacubaya'),
    @('B1027', 'Mazatlán Villa De Flores'),
    @('B1029', 'Miahuatlán De Porfirio Díaz'),
    @('B1030', 'Mixistlán De La Reforma'),
    @('B1034', 'Nejapa De Madero'),
    @('B1036', 'Oaxaca De Juárez'),
    @('B1037', 'Ocotlán De Morelos'),
    @('B1038', 'Pinotepa De Don Luis'),
    @('B1040', 'Putla Villa De Guerrero'),
    @('B1041', 'Reforma De Pineda'),
    @('B1043', 'Rojas De Cuauhtémoc'),
    @('B1048', 'San Agustín De Las Juntas'),
    @('B1069', 'San Antonino El Alto'),
    @('B1072', 'San Antonio De La Cal'),
    @('B1079', 'San Baltazar Yatzachi El Bajo'),
    @('B1095', 'San Dionisio Del Mar'),
    @('B1099', 'San Felipe Jalapa De Díaz'),
    @('B1106', 'San Francisco Del Mar'),
    @('B1132', 'San José Del Peñasco'),
    @('B1133', 'San José Del Progreso'),
    @('B1145', 'San Juan Bautista Lo De Soto'),
    @('B1158', 'San Juan De Los Cués'),
    @('B1159', 'San Juan Del Estado'),
    @('B1160', 'San Juan Del Río'),
    @('B1200', 'San Martín De Los Cansecos'),
    @('B1208', 'San Mateo Del Mar'),
    @('B1226', 'San Miguel Del Puerto'),
    @('B1227', 'San Miguel Del Río'),
    @('B1229', 'San Miguel El Grande'),
    @('B1254', 'San Pablo Villa De Mitla'),
    @('B1261', 'San Pedro El Alto'),
    @('B1287', 'San Pedro Y San Pablo Ayutla'),
    @('B1288', 'San Pedro Y San Pablo Teposcolula'),
    @('B1289', 'San Pedro Y San Pablo Tequixtepec'),
    @('B1308', 'Santa Ana Del Valle'),
    @('B1327', 'Santa Cruz De Bravo'),
    @('B1332', 'Santa Cruz Tacache De Mina'),
    @('B1339', 'Santa Inés De Zaragoza'),
    @('B1340', 'Santa Inés Del Monte'),
    @('B1342', 'Santa Lucía Del Camino'),
    @('B1356', 'Santa María Del Tule'),
    @('B1364', 'Santa María Jalapa Del Marqués'),
    @('B1366', 'Santa María La Asunción'),
    @('B1406', 'Santiago Del Río'),
    @('B1446', 'Santo Domingo De Morelos'),
    @('B1472', 'Sitio De Xitlapehua'),
    @('B1474', 'Tamazulápam Del Espíritu Santo'),
    @('B1475', 'Tanetze De Zaragoza'),
    @('B1477', 'Tataltepec De Valdés'),
    @('B1478', 'Teococuilco De Marcos Pérez'),
    @('B1479', 'Teotitlán De Flores Magón'),
    @('B1480', 'Teotitlán Del Valle'),
    @('B1482', 'Tepelmeme Villa De Morelos'),
    @('B1483', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B1484', 'Tlacolula De Matamoros'),
    @('B1486', 'Tlalixtac De Cabrera'),
    @('B1487', 'Totontepec Villa De Morelos'),
    @('B1491', 'Villa De Chilapa De Díaz'),
    @('B1492', 'Villa De Etla'),
    @('B1493', 'Villa De Tamazulápam Del Progreso'),
    @('B1494', 'Villa De Tututepec'),
    @('B1495', 'Villa De Zaachila'),
    @('B1498', 'Villa Sola De Vega'),
    @('B1499', 'Villa Talea De Castro'),
    @('B1500', 'Villa Tejúpam De La Unión'),
    @('B1503', 'Yutanduchi De Guerrero'),
    @('B1504', 'Zapotitlán Del Río'),
    @('B1507', 'Zimatlán De Álvarez'),
    @('B1534', 'Ayotoxco De Guerrero'),
    @('B1538', 'Chalchicomula De Sesma'),
    @('B1548', 'Chila De La Sal'),
    @('B1559', 'Cuapiaxtla De Madero'),
    @('B1563', 'Cuayuca De Andrade'),
    @('B1564', 'Cuetzalan Del Progreso'),
    @('B1580', 'Huehuetlán El Chico'),
    @('B1581', 'Huehuetlán El Grande'),
    @('B1586', 'Huitzilan De Serdán'),
    @('B1588', 'Ixcamilpa De Guerrero'),
    @('B1592', 'Izúcar De Matamoros'),
    @('B1603', 'Los Reyes De Juárez'),
    @('B1604', 'Mazapiltepec De Juárez'),
    @('B1617', 'Palmar De Bravo'),
    @('B1627', 'San Diego La Mesa Tochimiltzingo'),
    @('B1644', 'San Nicolás De Los Ranchos'),
    @('B1648', 'San Salvador El Seco'),
    @('B1649', 'San Salvador El Verde'),
    @('B1658', 'Tecali De Herrera'),
    @('B1666', 'Tepanco De López'),
    @('B1667', 'Tepango De Rodríguez'),
    @('B1668', 'Tepatlaxco De Hidalgo'),
    @('B1673', 'Tepexi De Rodríguez'),
    @('B1675', 'Tepeyahualco De Cuauhtémoc'),
    @('B1676', 'Tetela De Ocampo'),
    @('B1677', 'Teteles De Avila Castillo'),
    @('B1682', 'Tlacotepec De Benito Juárez'),
    @('B1694', 'Totoltepec De Guerrero'),
    @('B1696', 'Tuzamapan De Galeana'),
    @('B1700', 'Xayacatlán De Bravo'),
    @('B1706', 'Xochitlán De Vicente Suárez'),
    @('B1721', 'Amealco De Bonfil'),
    @('B1723', 'Cadereyta De Montes'),
    @('B1729', 'Jalpan De Serra'),
    @('B1730', 'Landa De Matamoros'),
    @('B1733', 'Pinal De Amoles'),
    @('B1736', 'San Juan Del Río'),
    @('B1751', 'Armadillo De Los Infante'),
    @('B1752', 'Axtla De Terrazas'),
    @('B1758', 'Ciudad Del Maíz'),
    @('B1769', 'Mexquitic De Carmona'),
    @('B1775', 'San Ciro De Acosta'),
    @('B1781', 'Santa María Del Río'),
    @('B1783', 'Soledad De Graciano Sánchez'),
    @('B1791', 'Tanquián De Escobedo'),
    @('B1795', 'Villa De Arista'),
    @('B1796', 'Villa De Arriaga'),
    @('B1797', 'Villa De Guadalupe'),
    @('B1798', 'Villa De La Paz'),
    @('B1799', 'Villa De Ramos'),
    @('B1800', 'Villa De Reyes'),
    @('B1862', 'Nacozari De García'),
    @('B1874', 'San Felipe De Jesús'),
    @('B1877', 'San Miguel De Horcasitas'),
    @('B1895', 'Jalpa De Méndez'),
    @('B1929', 'Soto La Marina'),
    @('B1937', 'Acuamanala De Miguel Hidalgo'),
    @('B1939', 'Amaxac De Guerrero'),
    @('B1940', 'Apetatitlán De Antonio Carvajal'),
    @('B1946', 'Contla De Juan Cuamatzi'),
    @('B1953', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1957', 'Mazatecochco De José María Morelos'),
    @('B1958', 'Muñoz De Domingo Arenas'),
    @('B1959', 'Nanacamilpa De Mariano Arista'),
    @('B1962', 'Papalotla De Xicohténcatl'),
    @('B1968', 'San Pablo Del Monte'),
    @('B1976', 'Tepetitla De Lardizábal'),
    @('B1979', 'Tetla De La Solidaridad'),
    @('B1991', 'Ziltlaltépec De Trinidad Sánchez Santos'),
    @('B2001', 'Alto Lucero De Gutiérrez Barrios'),
    @('B2005', 'Amatlán De Los Reyes'),
    @('B2017', 'Boca Del Río'),
    @('B2019', 'Camarón De Tejeda'),
    @('B2023', 'Castillo De Teayo'),
    @('B2025', 'Cazones De Herrera'),
    @('B2045', 'Cosamaloapan De Carpio'),
    @('B2046', 'Cosautlán De Carvajal'),
    @('B2063', 'Hueyapan De Ocampo'),
    @('B2064', 'Huiloapan De Cuauhtémoc'),
    @('B2065', 'Ignacio De La Llave'),
    @('B2069', 'Ixhuacán De Los Reyes'),
    @('B2070', 'Ixhuatlán De Madero'),
    @('B2071', 'Ixhuatlán Del Café'),
    @('B2072', 'Ixhuatlán Del Sureste'),
    @('B2084', 'Juchique De Ferrer'),
    @('B2089', 'Las Vigas De Ramírez'),
    @('B2090', 'Lerdo De Tejada'),
    @('B2096', 'Martínez De La Torre'),
    @('B2099', 'Medellín De Bravo'),
    @('B2103', 'Mixtla De Altamirano'),
    @('B2105', 'Nanchital De Lázaro Cárdenas Del Río'),
    @('B2115', 'Ozuluama De Mascareñas'),
    @('B2119', 'Paso De Ovejas'),
    @('B2120', 'Paso Del Macho'),
    @('B2124', 'Poza Rica De Hidalgo'),
    @('B2134', 'Sayula De Alemán'),
    @('B2138', 'Soledad De Doblado'),
    @('B2146', 'Tatahuicapan De Juárez'),
    @('B2168', 'Tlacotepec De Mejía'),
    @('B2183', 'Vega De Alatorre'),
    @('B2195', 'Zontecomatlán De López Y Fuentes'),
    @('B2196', 'Zozocolco De Hidalgo'),
    @('B2287', 'Cañitas De Felipe Pescador'),
    @('B2289', 'Concepción Del Oro'),
    @('B2291', 'El Plateado De Joaquín Amaro'),
    @('B2301', 'Jiménez Del Teul'),
    @('B2307', 'Mezquital Del Oro'),
    @('B2312', 'Moyahua De Estrada'),
    @('B2313', 'Nochistlán De Mejía'),
    @('B2314', 'Noria De Ángeles'),
    @('B2325', 'Teúl De González Ortega'),
    @('B2326', 'Tlaltenango De Sánchez Román'),
    @('B2328', 'Trinidad García De La Cadena'),
    @('B2331', 'Villa De Cos')
)

foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Remove trailing footnote rows (2339:2344) so dimension shrinks to A1:D2338
$ws.Rows("2339:2344").Delete()
